$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns C:F (runs, balls, fours, sixes) for rows 2-12,
# reflecting the updated activity figures.
$data = @{
    2  = @("3", "2", "0", "0")
    3  = @("7", "4", "0", "1")
    4  = @("7", "9", "1", "0")
    5  = @("81", "51", "8", "3")
    6  = @("39", "34", "4", "2")
    7  = @("4", "10", "1", "0")
    8  = @("13", "12", "1", "0")
    9  = @("16", "22", "1", "0")
    10 = @("23", "16", "2", "1")
    11 = @("1", "5", "0", "0")
    12 = @("36", "16", "3", "3")
}

$cols = @("C", "D", "E", "F")

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $cell = $ws.Range($cols[$i] + $row)
        $cell.NumberFormat = "@"
        $cell.Value = $values[$i]
    }
}
